# Optuna Attempt (go back with original)
# Update forecast values on "Forecast Comparison" sheet and derived totals on "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---
# Row : MyForecast(D), Inventory Coverage(H), Seasonality Index(L)
$updates = @(
    @{ Row = 2;  D = 94;  H = 7.4;  L = 1.09 },
    @{ Row = 3;  D = 96;  H = 6.27; L = 0.87 },
    @{ Row = 4;  D = 96;  H = 5.27; L = 1.18 },
    @{ Row = 5;  D = 95;  H = 4.32; L = 1.08 },
    @{ Row = 6;  D = 96;  H = 3.28; L = 0.82 },
    @{ Row = 7;  D = 94;  H = 2.31; L = 0.98 },
    @{ Row = 8;  D = 91;  H = 1.36; L = 1.01 },
    @{ Row = 9;  D = 92;  H = 0.36; L = 0.94 },
    @{ Row = 10; L = 0.89 },
    @{ Row = 11; L = 0.96 },
    @{ Row = 12; L = 1.14 },
    @{ Row = 13; L = 0.89 },
    @{ Row = 14; L = 0.96 },
    @{ Row = 15; L = 0.9 },
    @{ Row = 16; L = 0.96 },
    @{ Row = 17; L = 0.82 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) { $wsForecast.Range("D$r").Value = $u.D }
    if ($u.ContainsKey("H")) { $wsForecast.Range("H$r").Value = $u.H }
    if ($u.ContainsKey("L")) { $wsForecast.Range("L$r").Value = $u.L }
}

# --- Summary sheet ---
# These "Value" cells are stored as text (not numbers) in the workbook, so a
# leading apostrophe is used to force Excel to keep them as text instead of
# auto-converting the numeric-looking string to a number.
$wsSummary.Range("B9").Formula2  = "'1413"
$wsSummary.Range("B10").Formula2 = "'755"
$wsSummary.Range("B11").Formula2 = "'381"
$wsSummary.Range("B12").Formula2 = "'96"
